$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A was empty -- delete it so everything (B:H) shifts left to (A:G).
$ws.Columns("A").Delete()

# The table used to span B6:H7; after the shift it lives at A6:G7. Resizing
# the ListObject keeps the table/autoFilter ref (and header dataDxfId links)
# in sync with the new location.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A6:G7"))

# New bold, centered title style for the merged header row.
$hdr = $ws.Range("A2")
$hdr.Font.Bold = $true
$hdr.Font.Size = 14
$hdr.Font.Name = "Arial"
$hdr.Font.Color = 4473924
$hdr.HorizontalAlignment = -4108
$ws.Range("A2:G2").Merge()

# Widen/resize the data columns to their new layout.
$ws.Columns("B").ColumnWidth = 11.307291666666666
$ws.Columns("C").ColumnWidth = 19.022135416666668
$ws.Columns("D").ColumnWidth = 13.022135416666666
$ws.Columns("E").ColumnWidth = 20.022135416666668
$ws.Columns("F").ColumnWidth = 27.451822916666668

# Print at 70% scale instead of 100%.
$ws.PageSetup.Zoom = 70

# Move the remembered selection down one row.
$ws.Range("F11").Select() | Out-Null
